# Insert a new weekly record above current row 191 ("Perejil" sheet),
# shifting the existing rows 191-219 down to 192-220, then populate the
# newly inserted row with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 191; everything below shifts down one row.
$ws.Rows.Item(191).Insert()

# Fill in the values for the new row 191.
$ws.Range("A191").Value = 8
$ws.Range("B191").Value = "Terminal La Palmera de La Serena"
$ws.Range("C191").Value = "Coquimbo"
$ws.Range("D191").Value = 45034
$ws.Range("E191").Value = 4
$ws.Range("F191").Value = 100112044
$ws.Range("G191").Value = "Perejil"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 2400
$ws.Range("K191").Value = 2000
$ws.Range("L191").Value = 2500
$ws.Range("M191").Value = 2250
$ws.Range("N191").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O191").Value = "Provincia del Elquí"
$ws.Range("P191").Value = 1500
$ws.Range("Q191").Value = 1.5
$ws.Range("R191").Value = "Hortaliza"
